# Apply the periodic cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sun Jun 30 20:40:32 UTC 2024 with GitHub Actions".
#
# Price (column D) cells are plain text in the workbook (values like "1.00" or
# "62.244.24" are not real numbers), so any assignment that LOOKS like a number
# is prefixed with a leading apostrophe to force Excel to keep it as text and
# preserve the exact digits/trailing zeros (Excel would otherwise silently
# coerce "580.05" -> 580.05 or "1.00" -> 1, dropping the formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "62.303.96"
$ws.Range("E2").Value = "  +2.10%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.431.61"
$ws.Range("E3").Value = "  +1.41%  "

# Row 4 (TetherUSD)
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'580.05"
$ws.Range("E5").Value = "  +1.45%  "

# Row 6 (Solana)
$ws.Range("D6").Value = "'145.19"
$ws.Range("E6").Value = "  +2.88%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 (XRP)
$ws.Range("E8").Value = "  +0.28%  "

# Row 9 (Toncoin)
$ws.Range("E9").Value = "  -0.67%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value = "  +1.40%  "

# Row 11 (Cardano)
$ws.Range("E11").Value = "  +0.24%  "

# Row 12 (WrappedliquidstakedEther2.0)
$ws.Range("D12").Value = "4.018.35"
$ws.Range("E12").Value = "  +1.37%  "

# Row 13 (Avalanche)
$ws.Range("D13").Value = "'28.99"
$ws.Range("E13").Value = "  +4.07%  "

# Row 14 (TRON)
$ws.Range("E14").Value = "  -0.68%  "

# Row 15 (WrappedEther)
$ws.Range("D15").Value = "3.429.54"
$ws.Range("E15").Value = "  +1.27%  "

# Row 16 (ShibaInu)
$ws.Range("E16").Value = "  +0.56%  "

# Row 17 (WrappedBTC)
$ws.Range("D17").Value = "62.305.28"
$ws.Range("E17").Value = "  +1.96%  "

# Row 18 (Polkadot)
$ws.Range("E18").Value = "  +1.67%  "

# Row 19 (Chainlink)
$ws.Range("D19").Value = "'14.07"
$ws.Range("E19").Value = "  +2.95%  "

# Row 20 (Uniswap)
$ws.Range("E20").Value = "  +3.05%  "

# Row 21 (BitcoinCash)
$ws.Range("D21").Value = "'393.90"
$ws.Range("E21").Value = "  +2.44%  "

# Row 22 (Litecoin)
$ws.Range("D22").Value = "'74.98"
$ws.Range("E22").Value = "  -2.00%  "

# Row 23 (Polygon)
$ws.Range("E23").Value = "  +0.62%  "

# Row 24 (Dai)
$ws.Range("E24").Value = "  +0.04%  "

# Row 25 (PEPE)
$ws.Range("D25").Value = "'0.0000116"
$ws.Range("E25").Value = "  +1.13%  "

# Row 26 (WrappedeETH)
$ws.Range("D26").Value = "3.568.44"
$ws.Range("E26").Value = "  +1.30%  "

# Row 27 (Kaspa)
$ws.Range("D27").Value = "'0.187"
$ws.Range("E27").Value = "  +1.27%  "

# Row 28 (RenderToken)
$ws.Range("E28").Value = "  +4.65%  "

# Row 29 (Binance-PegBSC-USD)
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30 (InternetComputer(DFINITY))
$ws.Range("E30").Value = "  +1.30%  "

# Row 31 (PancakeSwap)
$ws.Range("E31").Value = "  +0.77%  "

# Row 32 (Fetch.AI)
$ws.Range("E32").Value = "  +2.81%  "

# Row 33 (USDe)
$ws.Range("E33").Value = "  +0.01%  "

# Row 34 (EthereumClassic)
$ws.Range("E34").Value = "  +1.53%  "

# Row 35 (NEARProtocol)
$ws.Range("D35").Value = "'5.34"
$ws.Range("E35").Value = "  +6.95%  "

# Row 36 (Aptos)
$ws.Range("D36").Value = "'7.01"
$ws.Range("E36").Value = "  +1.05%  "

# Row 37 (Monero)
$ws.Range("E37").Value = "  +1.28%  "

# Row 38 (ImmutableX)
$ws.Range("D38").Value = "'1.53"
$ws.Range("E38").Value = "  +5.02%  "

# Row 39 (RenzoRestakedETH)
$ws.Range("D39").Value = "3.464.37"
$ws.Range("E39").Value = "  +1.32%  "

# Row 40 (EnergySwap)
$ws.Range("D40").Value = "'28.95"
$ws.Range("E40").Value = "  +9.13%  "

# Row 41 (Hedera)
$ws.Range("D41").Value = "'0.0756"
$ws.Range("E41").Value = "  -1.21%  "

# Row 42 (Mantle)
$ws.Range("E42").Value = "  +1.74%  "

# Row 43 (Filecoin)
$ws.Range("E43").Value = "  +2.18%  "

# Row 44 (Stacks)
$ws.Range("E44").Value = "  +2.17%  "

# Row 45 (ONDO)
$ws.Range("D45").Value = "'1.17"
$ws.Range("E45").Value = "  +4.74%  "

# Row 46 (Maker)
$ws.Range("D46").Value = "2.516.46"
$ws.Range("E46").Value = "  +2.25%  "

# Row 47 (InjectiveProtocol)
$ws.Range("E47").Value = "  +1.42%  "

# Row 48 (Cosmos)
$ws.Range("D48").Value = "'6.70"
$ws.Range("E48").Value = "  +0.78%  "

# Row 49 (FirstDigitalUSD)
$ws.Range("D49").Value = "'1.00"
$ws.Range("E49").Value = "  -0.08%  "

# Row 50 (VeChain)
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "'2.15"
$ws.Range("E50").Value = "  +1.02%  "

# Row 51 (dogwifhat)
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0264"
$ws.Range("E51").Value = "  +0.85%  "
